$d = $word.ActiveDocument

# --- Paragraph 1: "Attended by: Angel, Georgi, Ilia, Mikaeil" -------------
# Strike out the whole paragraph mark (pPr/rPr/strike) and give each of the
# three names ("Angel", "Georgi", "Mikaeil") their own strike-through run,
# with spell-check proofErr markers bracketing "Georgi" and "Mikaeil" just
# like the target OOXML.
$p1xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:t>Attended by:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>Angel</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>Georgi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Ilia</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>Mikaeil</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# --- Paragraph 2: replaces "Planned activities:" text with the new line --
# "Nobody said anything." (split into two runs, matching the target), while
# keeping the existing _GoBack bookmark anchored to this paragraph.
$p2xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Nobody said anything</w:t></w:r><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML($p1xml)

$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML($p2xml)

# --- New paragraph holding "Planned activities:" --------------------------
# The old "Planned activities:" paragraph became "Nobody said anything."
# above, so re-create a fresh paragraph right after it for the original
# "Planned activities:" line (ahead of "Highlights:").
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Planned activities:"
